# EPBDS-11247: use Integer.valueOf instead of deprecated constructor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value  = "'= Integer.valueOf(20)"
$ws.Range("C34").Value = "'= Integer.valueOf(20)"
$ws.Range("C65").Value = "'= Integer.valueOf(20)"
$ws.Range("C93").Value = "'= Integer.valueOf(20)"

$ws.Range("D36").Value = "'= Integer.valueOf(200)"
$ws.Range("D67").Value = "'= Integer.valueOf(200)"
$ws.Range("D95").Value = "'= Integer.valueOf(200)"

$ws.Range("D107").Select()
